# Update countries & provincias Spain
# - Reorders a couple of shared-string country names (Nepal before
#   "Sudan del Sur"; "Etiopia" moved up next to "Tanzania"/"Madagascar")
#   and refreshes their COVID case numbers.
# - Bumps the "datos actualizados" timestamp.
# - Refreshes the numeric stats for Iran (row 13) and the block of
#   countries around rows 129-138.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Iran: updated totals
Set-Row 13 "Iran" 133521 1869 104072 22090 0 59 7359

# Nepal now appears before "Sudan del Sur" (row 129/130 swap), both with
# refreshed data for Nepal and unchanged data carried over for Sudan del Sur.
Set-Row 129 "Nepal" 584 68 70 511 0 0 3
Set-Row 130 "Sudan del Sur" 563 0 6 551 0 0 6

# "Etiopia" moves up to row 133 (right after Tanzania) with refreshed data;
# the countries that used to sit there (Madagascar, Republica de Africa
# Central, Congo, Reunion, Taiwan) each shift down one row, keeping their
# previous data.
Set-Row 133 "Etiopia" 494 61 151 338 0 0 5
Set-Row 134 "Madagascar" 488 40 138 348 0 0 2
Set-Row 135 "Republica de Africa Central" 479 0 18 461 0 0 0
Set-Row 136 "Congo" 469 0 137 316 0 0 16
Set-Row 137 "Reunion" 449 0 411 37 0 0 1
Set-Row 138 "Taiwan" 441 0 411 23 0 0 7

# Bump the "last updated" timestamp string (A1).
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 13:05"
